$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Fill in the "Actual Result" column (J) for rows 11-18 in order
$ws.Range("J11").Value = "User was successfully logged in and redirected to the Dashboard page."
$ws.Range("J12").Value = "System displayed a warning message “Invalid credentials” and prevented login."
$ws.Range("J13").Value = "System displayed a warning message “Invalid credentials” and prevented login."
$ws.Range("J14").Value = "System displayed a warning message “Invalid credentials” and prevented login."
$ws.Range("J15").Value = "User was successfully logged out and redirected to the login page."
$ws.Range("J16").Value = "New employee “WadhaAlgarni2” was successfully created and appeared in the Employee List."
$ws.Range("J17").Value = "Employee record for “WadhaAlgarni2” appeared in the Employee List with correct details."
$ws.Range("J18").Value = "System displayed validation messages for required fields and did not create a new record."

# Rows 19-22 were filled in reverse order (22, 21, 20, 19)
$ws.Range("J22").Value = "Personal details were successfully updated and saved in the system."
$ws.Range("J21").Value = "Profile image was successfully uploaded and updated on the user profile."
$ws.Range("J20").Value = "System displayed an error “Job Title already exists” and did not save the duplicate job."
$ws.Range("J19").Value = "Job “Software Engineer - QA” was successfully added and appeared in the Job Titles list."

# Fill in the "Status" column (K) with "Pass" for all rows
$ws.Range("K11").Value = "Pass"
$ws.Range("K12").Value = "Pass"
$ws.Range("K13").Value = "Pass"
$ws.Range("K14").Value = "Pass"
$ws.Range("K15").Value = "Pass"
$ws.Range("K16").Value = "Pass"
$ws.Range("K17").Value = "Pass"
$ws.Range("K18").Value = "Pass"
$ws.Range("K19").Value = "Pass"
$ws.Range("K20").Value = "Pass"
$ws.Range("K21").Value = "Pass"
$ws.Range("K22").Value = "Pass"
